$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.271.93'
$ws.Range('E2').Value = '  -2.14%  '
$ws.Range('D3').Value = '3.496.70'
$ws.Range('E3').Value = '  -2.42%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'596.53"
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('D6').Value = "'142.01"
$ws.Range('E6').Value = '  -2.48%  '
$ws.Range('D7').Value = '3.496.54'
$ws.Range('E7').Value = '  -2.36%  '
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').Value = "'0.502"
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('D10').Value = "'0.132"
$ws.Range('E10').Value = '  -3.22%  '
$ws.Range('D11').Value = "'7.77"
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('D12').Value = "'0.400"
$ws.Range('E12').Value = '  -3.92%  '
$ws.Range('D13').Value = '4.097.58'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').Value = "'0.0000196"
$ws.Range('E14').Value = '  -6.02%  '
$ws.Range('D15').Value = "'28.19"
$ws.Range('E15').Value = '  -6.16%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.516.86'
$ws.Range('E16').Value = '  -3.06%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = "'0.117"
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('D18').Value = '65.374.30'
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('D19').Value = "'10.97"
$ws.Range('E19').Value = '  -4.01%  '
$ws.Range('D20').Value = "'6.14"
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').Value = "'14.17"
$ws.Range('E21').Value = '  -5.47%  '
$ws.Range('D22').Value = "'412.52"
$ws.Range('E22').Value = '  -4.55%  '
$ws.Range('D23').Value = "'0.591"
$ws.Range('E23').Value = '  -4.67%  '
$ws.Range('D24').Value = "'77.20"
$ws.Range('E24').Value = '  -2.53%  '
$ws.Range('D25').Value = '3.645.81'
$ws.Range('E25').Value = '  -2.11%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = "'0.0000113"
$ws.Range('E27').Value = '  -6.24%  '
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').Value = "'7.61"
$ws.Range('E29').Value = '  -5.72%  '
$ws.Range('D30').Value = "'8.72"
$ws.Range('E30').Value = '  -6.29%  '
$ws.Range('D31').Value = "'0.997"
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').Value = '3.509.06'
$ws.Range('E32').Value = '  -1.92%  '
$ws.Range('D33').Value = "'0.151"
$ws.Range('E33').Value = '  -2.92%  '
$ws.Range('D34').Value = "'24.08"
$ws.Range('E34').Value = '  -5.36%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').Value = "'7.43"
$ws.Range('E36').Value = '  -5.27%  '
$ws.Range('D37').Value = "'1.27"
$ws.Range('E37').Value = '  -12.34%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = "'172.69"
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = "'5.21"
$ws.Range('E39').Value = '  -7.47%  '
$ws.Range('D40').Value = "'1.55"
$ws.Range('E40').Value = '  -10.12%  '
$ws.Range('D41').Value = "'0.0809"
$ws.Range('E41').Value = '  -5.26%  '
$ws.Range('D42').Value = "'5.01"
$ws.Range('E42').Value = '  -4.20%  '
$ws.Range('D43').Value = "'0.850"
$ws.Range('E43').Value = '  -4.83%  '
$ws.Range('D44').Value = "'45.06"
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('D45').Value = "'1.75"
$ws.Range('E45').Value = '  -9.26%  '
$ws.Range('D46').Value = "'1.00"
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').Value = "'2.33"
$ws.Range('E47').Value = '  -8.30%  '
$ws.Range('D48').Value = "'6.96"
$ws.Range('E48').Value = '  -3.58%  '
$ws.Range('D49').Value = "'22.74"
$ws.Range('E49').Value = '  -3.64%  '
$ws.Range('E50').Value = '  -10.56%  '
$ws.Range('D51').Value = '2.337.99'
$ws.Range('E51').Value = '  -1.79%  '
